# Applies the "changing excel data" commit to Book1.xlsx.
#
# Summary of the change (student roster swap + ripple into the dependent
# DegreePlan / Student / StudentTerm sheets):
#   - DegreePlan (sheet5): rows now reference the new student IDs
#     (S533990 / S531366 / S533710 / S533626 / S531383), two extra rows
#     (10 & 11) added for the 5th student.
#   - Student (sheet6): the 8 old S511111..S511118 rows are replaced by
#     5 new students (real names), the phone-number column is only kept
#     for the first of them (with an updated number), rows 7-9 removed.
#   - StudentTerm (sheet7): the "StudentId(PK)" header typo is corrected
#     to "StudentID(FK)", and rows 3-6 get the new student IDs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DegreePlan sheet
# ---------------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("DegreePlan")

$wsPlan.Cells.Item(2,1).Value  = 12
$wsPlan.Cells.Item(2,2).Value  = 2
$wsPlan.Cells.Item(2,3).Value  = "S533990"
$wsPlan.Cells.Item(2,4).Value  = "No Summer Off"
$wsPlan.Cells.Item(2,5).Value  = "No Summer Off"

$wsPlan.Cells.Item(3,1).Value  = 13
$wsPlan.Cells.Item(3,2).Value  = 2
$wsPlan.Cells.Item(3,3).Value  = "S533990"
$wsPlan.Cells.Item(3,4).Value  = "Summer Off"
$wsPlan.Cells.Item(3,5).Value  = "Summer Off"

$wsPlan.Cells.Item(4,1).Value  = 12
$wsPlan.Cells.Item(4,2).Value  = 2
$wsPlan.Cells.Item(4,3).Value  = "S531366"
$wsPlan.Cells.Item(4,4).Value  = "No Summer Off"
$wsPlan.Cells.Item(4,5).Value  = "No Summer Off"

$wsPlan.Cells.Item(5,1).Value  = 13
$wsPlan.Cells.Item(5,2).Value  = 2
$wsPlan.Cells.Item(5,3).Value  = "S531366"
$wsPlan.Cells.Item(5,4).Value  = "Summer Off"
$wsPlan.Cells.Item(5,5).Value  = "Summer Off"

$wsPlan.Cells.Item(6,1).Value  = 12
$wsPlan.Cells.Item(6,2).Value  = 2
$wsPlan.Cells.Item(6,3).Value  = "S533710"
$wsPlan.Cells.Item(6,4).Value  = "No Summer Off"
$wsPlan.Cells.Item(6,5).Value  = "No Summer Off"

$wsPlan.Cells.Item(7,1).Value  = 13
$wsPlan.Cells.Item(7,2).Value  = 2
$wsPlan.Cells.Item(7,3).Value  = "S533710"
$wsPlan.Cells.Item(7,4).Value  = "Summer Off"
$wsPlan.Cells.Item(7,5).Value  = "Summer Off"

$wsPlan.Cells.Item(8,1).Value  = 12
$wsPlan.Cells.Item(8,2).Value  = 2
$wsPlan.Cells.Item(8,3).Value  = "S533626"
$wsPlan.Cells.Item(8,4).Value  = "No Summer Off"
$wsPlan.Cells.Item(8,5).Value  = "No Summer Off"

$wsPlan.Cells.Item(9,1).Value  = 13
$wsPlan.Cells.Item(9,2).Value  = 2
$wsPlan.Cells.Item(9,3).Value  = "S533626"
$wsPlan.Cells.Item(9,4).Value  = "Summer Off"
$wsPlan.Cells.Item(9,5).Value  = "Summer Off"

# Two brand-new rows for the 5th student
$wsPlan.Cells.Item(10,1).Value = 12
$wsPlan.Cells.Item(10,2).Value = 2
$wsPlan.Cells.Item(10,3).Value = "S531383"
$wsPlan.Cells.Item(10,4).Value = "No Summer Off"
$wsPlan.Cells.Item(10,5).Value = "No Summer Off"

$wsPlan.Cells.Item(11,1).Value = 13
$wsPlan.Cells.Item(11,2).Value = 2
$wsPlan.Cells.Item(11,3).Value = "S531383"
$wsPlan.Cells.Item(11,4).Value = "Summer Off"
$wsPlan.Cells.Item(11,5).Value = "Summer Off"

$wsPlan.Activate()
$wsPlan.Range("C10").Select()

# ---------------------------------------------------------------------
# Student sheet
# ---------------------------------------------------------------------
$wsStudent = $wb.Worksheets.Item("Student")

$wsStudent.Cells.Item(2,1).Value = "S533990"
$wsStudent.Cells.Item(2,2).Value = "Hari Priya"
$wsStudent.Cells.Item(2,3).Value = "Jupally"
$wsStudent.Cells.Item(2,4).Value = 919569151

$wsStudent.Cells.Item(3,1).Value = "S531366"
$wsStudent.Cells.Item(3,2).Value = "Bharadwaj"
$wsStudent.Cells.Item(3,3).Value = "Dasari"

$wsStudent.Cells.Item(4,1).Value = "S533710"
$wsStudent.Cells.Item(4,2).Value = "Durga Susmitha"
$wsStudent.Cells.Item(4,3).Value = "Kotyada"

$wsStudent.Cells.Item(5,1).Value = "S533626"
$wsStudent.Cells.Item(5,2).Value = "Hyndavi"
$wsStudent.Cells.Item(5,3).Value = "Musipatla"

$wsStudent.Cells.Item(6,1).Value = "S531383"
$wsStudent.Cells.Item(6,2).Value = "Vyshnavi"
$wsStudent.Cells.Item(6,3).Value = "Yalamareddy"

# Phone numbers only survive for row 2 (with a new value); rows 3-6 lose
# their "D" cell entirely, and rows 7-9 (the old extra students) go away.
$wsStudent.Range("D3:D6").ClearContents()
$wsStudent.Rows("7:9").Delete()

$wsStudent.Activate()
$wsStudent.Range("A2:A6").Select()

# ---------------------------------------------------------------------
# StudentTerm sheet
# ---------------------------------------------------------------------
$wsTerm = $wb.Worksheets.Item("StudentTerm")

# Fix the "StudentId(PK)" typo -> "StudentID(FK)", matching the rest of
# the workbook's *(FK)" header convention.
$wsTerm.Cells.Item(1,2).Value = "StudentID(FK)"
$wsTerm.Cells.Item(1,3).Value = "TermID(FK)"

$wsTerm.Cells.Item(3,2).Value = "S531366"
$wsTerm.Cells.Item(4,2).Value = "S533710"
$wsTerm.Cells.Item(5,2).Value = "S533626"
$wsTerm.Cells.Item(6,2).Value = "S531383"

$wsTerm.Activate()
$wsTerm.Range("A2").Select()
